# Adds five new variable rows (tree_planting_hours_3, tree_planting_costs_3,
# maintaining_trees_3, mulch_h_3, tree_var_costs) to the bottom of the
# "Tabelle1" parameter table, mirroring the formatting of the existing rows,
# and updates the sheet selection to reflect where the user ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlShiftDown = -4121, xlFormatFromLeftOrAbove = -4161
$xlShiftDown = -4121
$xlFormatFromLeftOrAbove = -4161

# --- Row 45: tree_planting_hours_3 ---------------------------------------
$ws.Rows(45).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$ws.Range("A45").Value = "tree_planting_hours_3"
$ws.Range("B45").Value = 6
$ws.Range("C45").Value = "NA"
$ws.Range("D45").Value = 10
$ws.Range("E45").Value = "posnorm"
$ws.Range("F45").Value = "Planting hours"

# --- Row 46: tree_planting_costs_3 ----------------------------------------
$ws.Rows(46).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$ws.Range("A46").Value = "tree_planting_costs_3"
$ws.Range("B46").Value = 3700
$ws.Range("C46").Value = "NA"
$ws.Range("D46").Value = 4200
$ws.Range("E46").Value = "posnorm"
$ws.Range("F46").Value = "Price for trees"

# --- Row 47: maintaining_trees_3 ------------------------------------------
$ws.Rows(47).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$ws.Range("A47").Value = "maintaining_trees_3"
$ws.Range("B47").Value = 3
$ws.Range("C47").Value = "NA"
$ws.Range("D47").Value = 6
$ws.Range("E47").Value = "posnorm"
$ws.Range("F47").Value = "Cultivate trees"

# --- Row 48: mulch_h_3 -----------------------------------------------------
$ws.Rows(48).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$ws.Range("A48").Value = "mulch_h_3"
$ws.Range("B48").Value = 2
$ws.Range("C48").Value = "NA"
$ws.Range("D48").Value = 4
$ws.Range("E48").Value = "posnorm"
$ws.Range("F48").Value = "Mulching hours"

# --- Row 49: tree_var_costs -------------------------------------------------
$ws.Rows(49).Insert($xlShiftDown, $xlFormatFromLeftOrAbove)
$ws.Range("A49").Value = "tree_var_costs"
$ws.Range("B49").Value = 100
$ws.Range("C49").Value = "NA"
$ws.Range("D49").Value = 150
$ws.Range("E49").Value = "posnorm"
$ws.Range("F49").Value = "Other tree costs"

# Reflect the final cursor/selection position left in the workbook.
$ws.Range("A50").Select() | Out-Null
